$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shorten the "task" column labels (B2:B16).
# Order matters for shared-string append order: set H17 first so its new
# string slots in ahead of "Enter"/"Prepare"/"Model".
$ws.Range("H17").Value = "filename = : required for the function"

$ws.Range("B2").Value = "Enter"
$ws.Range("B3").Value = "Enter"
$ws.Range("B4").Value = "Prepare"
$ws.Range("B5").Value = "Prepare"
$ws.Range("B6").Value = "Prepare"
$ws.Range("B7").Value = "Prepare"
$ws.Range("B8").Value = "Prepare"
$ws.Range("B9").Value = "Prepare"
$ws.Range("B10").Value = "Prepare"
$ws.Range("B11").Value = "Prepare"
$ws.Range("B12").Value = "Prepare"
$ws.Range("B13").Value = "Prepare"
$ws.Range("B14").Value = "Prepare"
$ws.Range("B15").Value = "Prepare"
$ws.Range("B16").Value = "Model"

# Move the active selection on Sheet1 from K20 to B16.
[void]$ws.Range("B16").Select()

# Best-effort: nudge the saved window position (xWindow 4360 -> 5860).
# The host app may not persist window chrome through this COM surface, but
# this mirrors the actual user action of moving the Excel window.
$win = $excel.ActiveWindow
if ($win) {
    $win.Left = 5860
}
